$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading-percent results for the 380 kV case (rows 2-25 correspond to
# time steps 0-23; columns B,C,D,E,F,H,I,K hold the per-line loading values).
$arr = New-Object 'object[,]' 24,10

$arr[0,0] = 10.15175217277664
$arr[0,1] = 7.591137565417258
$arr[0,2] = 6.636188441706111
$arr[0,3] = 16.3150468111069
$arr[0,4] = 36.77115165190033
$arr[0,5] = 0
$arr[0,6] = 7.344005520526261
$arr[0,7] = 28.50331753538476
$arr[0,8] = 0
$arr[0,9] = 11.09120877063395

$arr[1,0] = 9.94799364222359
$arr[1,1] = 7.305578535318151
$arr[1,2] = 6.626819512963642
$arr[1,3] = 15.39923975154767
$arr[1,4] = 36.27289477116938
$arr[1,5] = 0
$arr[1,6] = 7.344005520526261
$arr[1,7] = 28.26890353977643
$arr[1,8] = 0
$arr[1,9] = 10.93826670383468

$arr[2,0] = 9.825094547100854
$arr[2,1] = 7.127311570613736
$arr[2,2] = 6.620946187909641
$arr[2,3] = 14.81373973783732
$arr[2,4] = 35.97052100670538
$arr[2,5] = 0
$arr[2,6] = 7.344005520526261
$arr[2,7] = 28.12803250588073
$arr[2,8] = 0
$arr[2,9] = 10.84749309914064

$arr[3,0] = 9.775654383845565
$arr[3,1] = 7.054060930337442
$arr[3,2] = 6.618521889491574
$arr[3,3] = 14.56958767599028
$arr[3,4] = 35.84831636770147
$arr[3,5] = 0
$arr[3,6] = 7.344005520526261
$arr[3,7] = 28.07143124193673
$arr[3,8] = 0
$arr[3,9] = 10.81134135927829

$arr[4,0] = 9.767486197857814
$arr[4,1] = 7.041865241696153
$arr[4,2] = 6.618117457203905
$arr[4,3] = 14.52871993212616
$arr[4,4] = 35.82808916676542
$arr[4,5] = 0
$arr[4,6] = 7.344005520526261
$arr[4,7] = 28.06208225300238
$arr[4,8] = 0
$arr[4,9] = 10.80539053999745

$arr[5,0] = 9.824425065198854
$arr[5,1] = 7.126325954982056
$arr[5,2] = 6.620913618670878
$arr[5,3] = 14.81046911098301
$arr[5,4] = 35.96886864477614
$arr[5,5] = 0
$arr[5,6] = 7.344005520526261
$arr[5,7] = 28.1272658575189
$arr[5,8] = 0
$arr[5,9] = 10.84700207977248

$arr[6,0] = 10.08108961371379
$arr[6,1] = 7.493360986721521
$arr[6,2] = 6.632982797120013
$arr[6,3] = 16.00423940336701
$arr[6,4] = 36.59868708436414
$arr[6,5] = 0
$arr[6,6] = 7.344005520526261
$arr[6,7] = 28.42187816365405
$arr[6,8] = 0
$arr[6,9] = 11.03785534076233

$arr[7,0] = 10.59815332256745
$arr[7,1] = 8.184654637037276
$arr[7,2] = 6.655715116606745
$arr[7,3] = 18.17758843761275
$arr[7,4] = 37.85662320908641
$arr[7,5] = 0
$arr[7,6] = 7.344005520526261
$arr[7,7] = 29.02226980716637
$arr[7,8] = 0
$arr[7,9] = 11.43464489911199

$arr[8,0] = 10.98145360590112
$arr[8,1] = 8.669085516849982
$arr[8,2] = 6.671878370267782
$arr[8,3] = 19.80365641746469
$arr[8,4] = 38.78754424259912
$arr[8,5] = 0
$arr[8,6] = 7.344005520526261
$arr[8,7] = 29.47493120445522
$arr[8,8] = 0
$arr[8,9] = 11.73669980367531

$arr[9,0] = 11.15556761358648
$arr[9,1] = 8.883309059300858
$arr[9,2] = 6.679121392432569
$arr[9,3] = 20.50222134586841
$arr[9,4] = 39.2109804315716
$arr[9,5] = 0
$arr[9,6] = 7.344005520526261
$arr[9,7] = 29.68286424505726
$arr[9,8] = 0
$arr[9,9] = 11.87571481866971

$arr[10,0] = 11.22138549805399
$arr[10,1] = 8.963471857356481
$arr[10,2] = 6.681848880390545
$arr[10,3] = 20.76087198558857
$arr[10,4] = 39.37119322328395
$arr[10,5] = 0
$arr[10,6] = 7.344005520526261
$arr[10,7] = 29.76184915611459
$arr[10,8] = 0
$arr[10,9] = 11.92853102948687

$arr[11,0] = 11.20721678304492
$arr[11,1] = 8.946251110841345
$arr[11,2] = 6.681262139194143
$arr[11,3] = 20.70542792865352
$arr[11,4] = 39.33669644997202
$arr[11,5] = 0
$arr[11,6] = 7.344005520526261
$arr[11,7] = 29.74482815217316
$arr[11,8] = 0
$arr[11,9] = 11.91714926070671

$arr[12,0] = 11.16098524004246
$arr[12,1] = 8.889923717012195
$arr[12,2] = 6.679346088097268
$arr[12,3] = 20.52361834643389
$arr[12,4] = 39.22416462248172
$arr[12,5] = 0
$arr[12,6] = 7.344005520526261
$arr[12,7] = 29.68935772435225
$arr[12,8] = 0
$arr[12,9] = 11.88005687099317

$arr[13,0] = 11.1326497901735
$arr[13,1] = 8.855294605755059
$arr[13,2] = 6.678170475089659
$arr[13,3] = 20.41148992871482
$arr[13,4] = 39.1552145518269
$arr[13,5] = 0
$arr[13,6] = 7.344005520526261
$arr[13,7] = 29.65541108466839
$arr[13,8] = 0
$arr[13,9] = 11.85735772757747

$arr[14,0] = 10.97006292632463
$arr[14,1] = 8.654955173816955
$arr[14,2] = 6.671402855880358
$arr[14,3] = 19.75717871070619
$arr[14,4] = 38.7598598844659
$arr[14,5] = 0
$arr[14,6] = 7.344005520526261
$arr[14,7] = 29.46137903951747
$arr[14,8] = 0
$arr[14,9] = 11.72764231236643

$arr[15,0] = 10.8701991545724
$arr[15,1] = 8.530422552459211
$arr[15,2] = 6.667223388640825
$arr[15,3] = 19.34526785555446
$arr[15,4] = 38.51722005245452
$arr[15,5] = 0
$arr[15,6] = 7.344005520526261
$arr[15,7] = 29.34283121471909
$arr[15,8] = 0
$arr[15,9] = 11.64843719397132

$arr[16,0] = 10.81274099451772
$arr[16,1] = 8.458219986730962
$arr[16,2] = 6.664809112956458
$arr[16,3] = 19.10447308472854
$arr[16,4] = 38.37766329313547
$arr[16,5] = 0
$arr[16,6] = 7.344005520526261
$arr[16,7] = 29.27483801834716
$arr[16,8] = 0
$arr[16,9] = 11.60303528837532

$arr[17,0] = 10.7932858390271
$arr[17,1] = 8.433677184982223
$arr[17,2] = 6.663989886885369
$arr[17,3] = 19.02227757645135
$arr[17,4] = 38.33041624344048
$arr[17,5] = 0
$arr[17,6] = 7.344005520526261
$arr[17,7] = 29.25185108178416
$arr[17,8] = 0
$arr[17,9] = 11.58769132416084

$arr[18,0] = 10.88083239359311
$arr[18,1] = 8.543739316717179
$arr[18,2] = 6.667669368414418
$arr[18,3] = 19.3895172117074
$arr[18,4] = 38.5430500685966
$arr[18,5] = 0
$arr[18,6] = 7.344005520526261
$arr[18,7] = 29.35543124451828
$arr[18,8] = 0
$arr[18,9] = 11.65685311165559

$arr[19,0] = 11.17456830442773
$arr[19,1] = 8.906495016207822
$arr[19,2] = 6.67990928975514
$arr[19,3] = 20.57717954842848
$arr[19,4] = 39.25722256270539
$arr[19,5] = 0
$arr[19,6] = 7.344005520526261
$arr[19,7] = 29.70564443693215
$arr[19,8] = 0
$arr[19,9] = 11.89094752541211

$arr[20,0] = 11.36583591518838
$arr[20,1] = 9.137958790976116
$arr[20,2] = 6.687819903982555
$arr[20,3] = 21.31912438952083
$arr[20,4] = 39.72314191584233
$arr[20,5] = 0
$arr[20,6] = 7.344005520526261
$arr[20,7] = 29.93593717192336
$arr[20,8] = 0
$arr[20,9] = 12.04493178158788

$arr[21,0] = 11.26384241929145
$arr[21,1] = 9.014958637815997
$arr[21,2] = 6.683605832079004
$arr[21,3] = 20.92625772293321
$arr[21,4] = 39.47458940734364
$arr[21,5] = 0
$arr[21,6] = 7.344005520526261
$arr[21,7] = 29.81291158756361
$arr[21,8] = 0
$arr[21,9] = 11.9626754450662

$arr[22,0] = 10.87602524075135
$arr[22,1] = 8.537720692802306
$arr[22,2] = 6.6674677766191
$arr[22,3] = 19.36952448323956
$arr[22,4] = 38.53137249501633
$arr[22,5] = 0
$arr[22,6] = 7.344005520526261
$arr[22,7] = 29.34973426441904
$arr[22,8] = 0
$arr[22,9] = 11.65304785295698

$arr[23,0] = 10.45732393523929
$arr[23,1] = 8.001364724351708
$arr[23,2] = 6.649662225365116
$arr[23,3] = 17.59255190232298
$arr[23,4] = 37.51461691874105
$arr[23,5] = 0
$arr[23,6] = 7.344005520526261
$arr[23,7] = 28.85765071692642
$arr[23,8] = 0
$arr[23,9] = 11.32521958770316

$ws.Range("B2:K25").Value = $arr
